$wb = $excel.ActiveWorkbook

# Sheet "OFF": update Road (row 3) target depth stats
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 233
$wsOff.Range("C3").Value = 171
$wsOff.Range("D3").Value = 54
$wsOff.Range("E3").Value = 37
$wsOff.Range("F3").Value = 5

# Sheet "DEF": update Road (row 3) target depth stats
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 381
$wsDef.Range("C3").Value = 283
$wsDef.Range("D3").Value = 95
$wsDef.Range("E3").Value = 46
